# Generate Report for Handoff
# Updates the handoff package GUID/file names, the zh-cn/de-de xliff
# filenames, and the associated handoff timestamps across the
# "Overview", "zh-cn" and "de-de" worksheets, keeping the hyperlink
# targets intact but refreshing their displayed text.

$wb = $excel.ActiveWorkbook

$oldGuid = "5ea0edb1-b92b-405b-a15b-0da865a2543a"
$newGuid = "24ef0adf-cf64-4a4c-a2e1-7657233ee48c"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

$oldZhXlf = "$oldGuid.4d8fbeef1dc67cb9af53c5a8a7c9347f785a0323.zh-cn.xlf"
$newZhXlf = "$newGuid.3caf69941f983abca11a7d60feb401da9f23f5fa.zh-cn.xlf"

$oldDeXlf = "$oldGuid.4d8fbeef1dc67cb9af53c5a8a7c9347f785a0323.de-de.xlf"
$newDeXlf = "$newGuid.3caf69941f983abca11a7d60feb401da9f23f5fa.de-de.xlf"

$oldOverviewDate = "2016-03-24 17:15:00"
$newOverviewDate = "2016-03-24 17:15:46"

$oldZhDate = "2016-03-24 17:14:53"
$newZhDate = "2016-03-24 17:15:41"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- Sheet "Overview" ---
# A2: handoff package markdown file name (hyperlink display + value)
# D2: latest handoff datetime (shared with de-de!E2)
$ovA2 = $wsOverview.Range("A2")

$ovA2.Value = $newMd
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($ovA2, "https://github.com/OpenLocalizationTest/oltest/blob/31e6693a6c3b4ae13b4942d4638e7fb5f00dc3c3/e2e/$oldMd", [System.Type]::Missing, [System.Type]::Missing, $newMd)

$wsOverview.Range("D2").Value = $newOverviewDate

# --- Sheet "zh-cn" ---
# A2: handoff package markdown file name
# D2: zh-cn xliff file name
# E2: zh-cn latest handoff datetime
$zhA2 = $wsZh.Range("A2")
$zhD2 = $wsZh.Range("D2")

$zhA2.Value = $newMd
$zhD2.Value = $newZhXlf
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($zhA2, "https://github.com/OpenLocalizationTest/oltest/blob/31e6693a6c3b4ae13b4942d4638e7fb5f00dc3c3/e2e/$oldMd", [System.Type]::Missing, [System.Type]::Missing, $newMd)
$wsZh.Hyperlinks.Add($zhD2, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/086a8fa3f2207de69a3c4442480517fd4bf57e38/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhXlf", [System.Type]::Missing, [System.Type]::Missing, $newZhXlf)

$wsZh.Range("E2").Value = $newZhDate

# --- Sheet "de-de" ---
# A2: handoff package markdown file name
# D2: de-de xliff file name
# E2: latest handoff datetime (shared with Overview!D2)
$deA2 = $wsDe.Range("A2")
$deD2 = $wsDe.Range("D2")

$deA2.Value = $newMd
$deD2.Value = $newDeXlf
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($deA2, "https://github.com/OpenLocalizationTest/oltest/blob/31e6693a6c3b4ae13b4942d4638e7fb5f00dc3c3/e2e/$oldMd", [System.Type]::Missing, [System.Type]::Missing, $newMd)
$wsDe.Hyperlinks.Add($deD2, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb9cb302d8b1b49e52737c8cbd621b902649575b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeXlf", [System.Type]::Missing, [System.Type]::Missing, $newDeXlf)

$wsDe.Range("E2").Value = $newOverviewDate
